$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.243.15"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "3.316.80"
$ws.Range("E3").Value = "  +6.05%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'601.61"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "'144.37"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.312.74"
$ws.Range("E8").Value = "  +6.24%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").Value = "'0.477"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'35.19"
$ws.Range("E14").Value = "  +3.25%  "
$ws.Range("D15").Value = "3.854.40"
$ws.Range("E15").Value = "  +5.89%  "
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "3.312.27"
$ws.Range("E17").Value = "  +5.79%  "
$ws.Range("D18").Value = "64.297.96"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'6.94"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").Value = "'486.67"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  +6.87%  "
$ws.Range("E23").Value = "  +5.83%  "
$ws.Range("D24").Value = "'13.67"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").Value = "'85.21"
$ws.Range("E25").Value = "  -2.79%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("D28").Value = "'8.36"
$ws.Range("E28").Value = "  +4.34%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("D32").Value = "'28.54"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").Value = "'2.60"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("D37").Value = "'53.31"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("E38").Value = "  +4.40%  "
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").Value = "'431.18"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").Value = "'2.81"
$ws.Range("E41").Value = "  +4.26%  "
$ws.Range("D42").Value = "'8.51"
$ws.Range("D43").Value = "3.024.85"
$ws.Range("E43").Value = "  +5.53%  "
$ws.Range("E44").Value = "  -5.06%  "
$ws.Range("E45").Value = "  +5.41%  "
$ws.Range("E46").Value = "  +6.82%  "
$ws.Range("D47").Value = "'26.49"
$ws.Range("E47").Value = "  +4.10%  "
$ws.Range("D48").Value = "'2.36"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D50").Value = "'0.116"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("D51").Value = "'35.14"
$ws.Range("E51").Value = "  +13.74%  "
